$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows - Day 0 row (row 2)
$ws.Cells.Item(2, 3).Value = 22.091571807861328
$ws.Cells.Item(2, 5).Value = 2081.0

# Day 1 row (row 4)
$ws.Cells.Item(4, 3).Value = 4.609101295471191
$ws.Cells.Item(4, 5).Value = 355.0

# Day 2 row (row 5)
$ws.Cells.Item(5, 3).Value = 3.637298822402954
$ws.Cells.Item(5, 5).Value = 374.0

# Day 3 row (row 6)
$ws.Cells.Item(6, 3).Value = 3.413088321685791
$ws.Cells.Item(6, 5).Value = 327.0

# Day 3 (dup) row (row 7)
$ws.Cells.Item(7, 3).Value = 3.637298822402954
$ws.Cells.Item(7, 5).Value = 374.0

# New row 8 - Day 4
$ws.Cells.Item(8, 1).Value = "Day 4"
$ws.Cells.Item(8, 2).Value = 4.026899814605713
$ws.Cells.Item(8, 3).Value = 3.18915057182312
$ws.Cells.Item(8, 4).Value = 549.0443725585938
$ws.Cells.Item(8, 5).Value = 311.0

# New row 9 - Day 5
$ws.Cells.Item(9, 1).Value = "Day 5"
$ws.Cells.Item(9, 2).Value = 3.6586999893188477
$ws.Cells.Item(9, 3).Value = 3.6056344509124756
$ws.Cells.Item(9, 4).Value = 659.986328125
$ws.Cells.Item(9, 5).Value = 340.0
